$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Chapter 3 & 4" Spanish translation in row 102 (column C):
# "Capítulo 3 & 4" -> "Capítulo 3 Y 4"
$ws.Range("C102").Value = "Capítulo 3 Y 4"

# Add the missing "Loop" translation as a new row 103
$ws.Range("A103").Value = "Loop"
$ws.Range("B103").Value = "Loop"
$ws.Range("C103").Value = "Bucle"

# Update selection to reflect where the user ended up after editing
$ws.Range("C106").Select() | Out-Null
